# Agrega la columna "diseño" (Fecha inicio / Fecha fin) a la fila 33
# de la tabla de pedido en la hoja de Visualizacion (Hoja1), marcando
# el item como terminado ("si") y con las fechas correspondientes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
if (-not $ws) { $ws = $wb.ActiveSheet }

# Terminado: de "no" a "si"
$ws.Range("C33").Value = "si"

# Fecha inicio / Fecha fin (copiando el formato de fecha ya usado en la tabla)
$ws.Range("D29:E29").Copy($ws.Range("D33:E33"))

# Deja la seleccion activa en D33:E33, como en el archivo final
$ws.Range("D33:E33").Select()
